# The source data now comes from a CSV ingest instead of the old JSON
# ingest, and the CSV no longer contains a "Bahrain" application-country
# row for 2014. Remove that single data row (and let every following row
# shift up), mirroring the removal of the now-unused "Bahrain" shared
# string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$found = $ws.Cells.Find("Bahrain")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}
